$wb = $excel.ActiveWorkbook

# Mapping of row -> new value for column F, applied to both the
# "展览" and "全部类型" worksheets (they mirror the same data).
$updates = @{
    4  = 352
    6  = 1853
    9  = 190
    13 = 4414
    15 = 327
    16 = 1207
    17 = 515
    19 = 781
    21 = 411
    23 = 207
    24 = 17
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
